$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F with people responsible for comparing algorithms (site compare, decision tree)
$ws.Range("F4").Value = "גל"
$ws.Range("F5").Value = "גל"
$ws.Range("F6").Value = "טלנ"
$ws.Range("F7").Value = "טלנ"
$ws.Range("F8").Value = "אנדריי"
$ws.Range("F9").Value = "אנדריי"

# Update selection to match the final cursor position
$ws.Range("F9").Select() | Out-Null
